$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1582184.8
$ws.Range("J17").Value = 1624939.9
$ws.Range("L17").Value = 4874819.699999999
$ws.Range("N17").Value = -4875155.699999999
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2500
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2850
$ws.Range("H113").Value = 13845.143
$ws.Range("I113").Value = 13845.143
$ws.Range("K113").Value = 13845.143
$ws.Range("M113").Value = -10591.143
$ws.Range("H116").Value = 3423058.2
$ws.Range("I116").Value = 4448806
$ws.Range("K116").Value = 4448806
$ws.Range("M116").Value = -4445364
$ws.Range("H127").Value = 881.8
$ws.Range("I127").Value = 800.7143
$ws.Range("K127").Value = 2402.1429
$ws.Range("M127").Value = 2557.8571
$ws.Range("H131").Value = 2164.5715
$ws.Range("I131").Value = 1192.8334
$ws.Range("J131").Value = 7995
$ws.Range("K131").Value = 3578.5002
$ws.Range("L131").Value = 23985
$ws.Range("M131").Value = 1461.4998
$ws.Range("N131").Value = -34065
$ws.Range("H137").Value = 83390.16
$ws.Range("I137").Value = 8714.143
$ws.Range("J137").Value = 104299.44
$ws.Range("K137").Value = 26142.429
$ws.Range("L137").Value = 312898.32
$ws.Range("M137").Value = -23592.429
$ws.Range("N137").Value = -317998.32
$ws.Range("H138").Value = 6011.5415
$ws.Range("I138").Value = 7329
$ws.Range("J138").Value = 5572.3887
$ws.Range("K138").Value = 21987
$ws.Range("L138").Value = 16717.1661
$ws.Range("M138").Value = -16847
$ws.Range("N138").Value = -26997.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3211.1316
$ws.Range("I32").Value = 2045.9854
$ws.Range("K32").Value = 2045.9854
$ws.Range("M32").Value = -1758.9854
$ws.Range("H132").Value = 3082212
$ws.Range("I132").Value = 1837.9333
$ws.Range("J132").Value = 10190768
$ws.Range("K132").Value = 5513.7999
$ws.Range("L132").Value = 30572304
$ws.Range("M132").Value = -2983.7999
$ws.Range("N132").Value = -30577364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 13486.223
$ws.Range("I22").Value = 16982
$ws.Range("J22").Value = 1251
$ws.Range("K22").Value = 16982
$ws.Range("L22").Value = 1251
$ws.Range("M22").Value = -16809
$ws.Range("N22").Value = -1597
$ws.Range("H134").Value = 10221.82
$ws.Range("I134").Value = 6703.9756
$ws.Range("K134").Value = 20111.9268
$ws.Range("M134").Value = -17576.9268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10935.1045
$ws.Range("I31").Value = 3936.4348
$ws.Range("J31").Value = 14593.5
$ws.Range("K31").Value = 3936.4348
$ws.Range("L31").Value = 14593.5
$ws.Range("M31").Value = -3641.4348
$ws.Range("N31").Value = -15183.5
$ws.Range("H34").Value = 10935.1045
$ws.Range("I34").Value = 3936.4348
$ws.Range("J34").Value = 14593.5
$ws.Range("K34").Value = 3936.4348
$ws.Range("L34").Value = 14593.5
$ws.Range("M34").Value = -3734.4348
$ws.Range("N34").Value = -14997.5
$ws.Range("H122").Value = 2875.2222
$ws.Range("I122").Value = 1123.8422
$ws.Range("K122").Value = 3371.5266
$ws.Range("M122").Value = -921.5266000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3806.0264
$ws.Range("I68").Value = 2259.4
$ws.Range("J68").Value = 4040.3635
$ws.Range("K68").Value = 6778.200000000001
$ws.Range("L68").Value = 12121.0905
$ws.Range("M68").Value = -5967.200000000001
$ws.Range("N68").Value = -13743.0905
$ws.Range("H71").Value = 3806.0264
$ws.Range("I71").Value = 2259.4
$ws.Range("J71").Value = 4040.3635
$ws.Range("K71").Value = 20334.6
$ws.Range("L71").Value = 36363.2715
$ws.Range("M71").Value = -16278.6
$ws.Range("N71").Value = -44475.2715
$ws.Range("H131").Value = 1456.35
$ws.Range("J131").Value = 1499.4517
$ws.Range("L131").Value = 4498.355100000001
$ws.Range("N131").Value = -14578.3551

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8662.883
$ws.Range("I70").Value = 6888.1763
$ws.Range("J70").Value = 10437.588
$ws.Range("K70").Value = 6888.1763
$ws.Range("L70").Value = 10437.588
$ws.Range("M70").Value = -6618.1763
$ws.Range("N70").Value = -10977.588
$ws.Range("H73").Value = 8662.883
$ws.Range("I73").Value = 6888.1763
$ws.Range("J73").Value = 10437.588
$ws.Range("K73").Value = 6888.1763
$ws.Range("L73").Value = 10437.588
$ws.Range("M73").Value = -5952.1763
$ws.Range("N73").Value = -12309.588

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4103.0625
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 4899.923
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 4899.923
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -5489.923
$ws.Range("H27").Value = 4103.0625
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 4899.923
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 4899.923
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -5113.923
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H40").Value = 21889.4
$ws.Range("I40").Value = 14862
$ws.Range("J40").Value = 49999
$ws.Range("K40").Value = 14862
$ws.Range("L40").Value = 49999
$ws.Range("M40").Value = -14726
$ws.Range("N40").Value = -50271
$ws.Range("H46").Value = 1256.5333
$ws.Range("I46").Value = 987.3333
$ws.Range("J46").Value = 2333.3333
$ws.Range("K46").Value = 987.3333
$ws.Range("L46").Value = 2333.3333
$ws.Range("M46").Value = -799.3333
$ws.Range("N46").Value = -2709.3333
$ws.Range("H55").Value = 198.44444
$ws.Range("I55").Value = 260
$ws.Range("K55").Value = 260
$ws.Range("M55").Value = -87
$ws.Range("H136").Value = 1069163.1
$ws.Range("I136").Value = 17002.666
$ws.Range("J136").Value = 2647404
$ws.Range("K136").Value = 51007.99800000001
$ws.Range("L136").Value = 7942212
$ws.Range("M136").Value = -48457.99800000001
$ws.Range("N136").Value = -7947312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5860449
$ws.Range("I132").Value = 1995
$ws.Range("J132").Value = 9766085
$ws.Range("K132").Value = 5985
$ws.Range("L132").Value = 29298255
$ws.Range("M132").Value = -3455
$ws.Range("N132").Value = -29303315

Write-Output "done"